$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("R1")
$ws1.Range("G2").Value = "3924:15:59"
$ws1.Range("G3").Value = "63:48:37"

$ws2 = $wb.Worksheets.Item("R2")
$ws2.Range("G2").Value = "12105:39:40"
$ws2.Range("G3").Value = "3235:23:09"
$ws2.Range("G4").Value = "473:34:43"

$ws4 = $wb.Worksheets.Item("R4")
$ws4.Range("G2").Value = "2951:29:29"
$ws4.Range("G3").Value = "178:41:44"

$ws5 = $wb.Worksheets.Item("R5")
$ws5.Range("G2").Value = "425:28:28"

$ws6 = $wb.Worksheets.Item("R6")
$ws6.Range("G2").Value = "66:00:46"
